# Updates cryptos list values (Price and Volume(1h) columns) to latest scrape.
# Generated from the authoritative diff of xl/worksheets/sheet1.xml.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "33.619.02"
$ws.Range("E2").Value = "  -1.02%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.766.18"
$ws.Range("E3").Value = "  -0.93%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "222.88"
$ws.Range("E5").Value = "  +0.59%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.543"
$ws.Range("E6").Value = "  -1.41%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.67"
$ws.Range("E8").Value = "  +0.69%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.288"
$ws.Range("E9").Value = "  +0.88%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0684"
$ws.Range("E10").Value = "  -3.55%  "

# Row 11
$ws.Range("E11").Value = "  +1.38%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.021.78"
$ws.Range("E12").Value = "  -0.87%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.00"
$ws.Range("E13").Value = "  +4.65%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.753.09"
$ws.Range("E14").Value = "  -1.73%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.676.58"
$ws.Range("E15").Value = "  -0.82%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.606"
$ws.Range("E16").Value = "  -3.13%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.10"
$ws.Range("E17").Value = "  -2.45%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.31"
$ws.Range("E18").Value = "  -2.29%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.69"
$ws.Range("E20").Value = "  -3.16%  "

# Row 21
$ws.Range("E21").Value = "  +0.11%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.50"
$ws.Range("E22").Value = "  -1.61%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.99"
$ws.Range("E23").Value = "  -1.76%  "

# Row 24
$ws.Range("E24").Value = "  -2.99%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.96"
$ws.Range("E25").Value = "  +0.92%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.04"
$ws.Range("E26").Value = "  -1.87%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.98"
$ws.Range("E27").Value = "  +0.06%  "

# Row 28
$ws.Range("E28").Value = "  -0.54%  "

# Row 29
$ws.Range("E29").Value = "  +0.20%  "

# Row 30
$ws.Range("E30").Value = "  +1.19%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0509"
$ws.Range("E31").Value = "  -1.97%  "

# Row 32
$ws.Range("E32").Value = "  -2.71%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.47"
$ws.Range("E33").Value = "  -0.45%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.78"
$ws.Range("E34").Value = "  -1.66%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.374.42"
$ws.Range("E35").Value = "  -2.20%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.642"
$ws.Range("E36").Value = "  +0.64%  "

# Row 37
$ws.Range("E37").Value = "  -2.60%  "

# Row 38
$ws.Range("E38").Value = "  -1.51%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.20"
$ws.Range("E39").Value = "  +4.89%  "

# Row 40
$ws.Range("E40").Value = "  +0.86%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "77.38"
$ws.Range("E41").Value = "  -2.55%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.65"
$ws.Range("E42").Value = "  -1.92%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.898"
$ws.Range("E43").Value = "  -3.80%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.39"
$ws.Range("E44").Value = "  +13.40%  "

# Row 45
$ws.Range("E45").Value = "  +4.37%  "

# Row 46
$ws.Range("E46").Value = "  +14.19%  "

# Row 47
$ws.Range("E47").Value = "  +1.05%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "106.47"
$ws.Range("E48").Value = "  +1.08%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.78"
$ws.Range("E49").Value = "  -2.47%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.921.16"
$ws.Range("E50").Value = "  -0.61%  "

# Row 51
$ws.Range("E51").Value = "  +0.28%  "
